$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the date column (A) formatting down to row 301 by copying the
# style of the last existing data row (A269) over the new rows.
$srcDate = $ws.Range("A269")
$dstDate = $ws.Range("A270:A301")
$srcDate.Copy($dstDate)

# New row data: row, date-serial, B, C, D
$rows = @(
    @(270, 44344, 0, 0, 0),
    @(271, 44345, 0, 0, 0),
    @(272, 44346, 0, 0, 0),
    @(273, 44347, 0, 0, 0),
    @(274, 44348, 0, 0, 0),
    @(275, 44349, 0, 0, 0),
    @(276, 44350, 1, 1, 62.34413965087282),
    @(277, 44351, 0, 1, 62.34413965087282),
    @(278, 44352, 0, 1, 62.34413965087282),
    @(279, 44353, 0, 1, 62.34413965087282),
    @(280, 44354, 0, 1, 62.34413965087282),
    @(281, 44355, 0, 1, 62.34413965087282),
    @(282, 44356, 0, 1, 62.34413965087282),
    @(283, 44357, 0, 0, 0),
    @(284, 44358, 0, 0, 0),
    @(285, 44359, 0, 0, 0),
    @(286, 44360, 0, 0, 0),
    @(287, 44361, 0, 0, 0),
    @(288, 44362, 0, 0, 0),
    @(289, 44363, 0, 0, 0),
    @(290, 44364, 0, 0, 0),
    @(291, 44365, 0, 0, 0),
    @(292, 44366, 0, 0, 0),
    @(293, 44367, 0, 0, 0),
    @(294, 44368, 0, 0, 0),
    @(295, 44369, 0, 0, 0),
    @(296, 44370, 0, 0, 0),
    @(297, 44371, 1, 1, 62.34413965087282),
    @(298, 44372, 0, 1, 62.34413965087282),
    @(299, 44373, 0, 1, 62.34413965087282),
    @(300, 44374, 0, 1, 62.34413965087282),
    @(301, 44375, 0, 1, 62.34413965087282)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}
